$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row with key "system.common.edit" below the existing last row (B8)
$ws.Range("B9").Value = "system.common.edit"

# Move the active selection to the newly added cell, matching the
# author's final cursor position in the edited workbook.
$ws.Range("B9").Select() | Out-Null
